$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6928
$ws1.Range("F4").Value = 48
$ws1.Range("F5").Value = 452
$ws1.Range("F7").Value = 6750
$ws1.Range("F8").Value = 68
$ws1.Range("F15").Value = 0
$ws1.Range("F17").Value = 48
$ws1.Range("F18").Value = 25
$ws1.Range("F19").Value = 12
$ws1.Range("F20").Value = 5121
$ws1.Range("F22").Value = 149
$ws1.Range("F23").Value = 540
$ws1.Range("F24").Value = 213
$ws1.Range("F25").Value = 216

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6928
$ws4.Range("F4").Value = 48
$ws4.Range("F7").Value = 6750
$ws4.Range("F8").Value = 68
$ws4.Range("F9").Value = 200
$ws4.Range("F11").Value = 19
$ws4.Range("F12").Value = 107
$ws4.Range("F13").Value = 405
$ws4.Range("F17").Value = 48
$ws4.Range("F19").Value = 12
$ws4.Range("F21").Value = 5121
$ws4.Range("F23").Value = 112
$ws4.Range("F24").Value = 149
$ws4.Range("F25").Value = 540
$ws4.Range("F26").Value = 213
$ws4.Range("F27").Value = 216
